$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match data (columns F:V) between row 134 and row 138.
#    (A:E stay the same for both rows — only the match/odds data moved.)
# ---------------------------------------------------------------------------
$row134 = $ws.Range("F134:V134").Value()
$row138 = $ws.Range("F138:V138").Value()

$ws.Range("F134:V134").Value = $row138
$ws.Range("F138:V138").Value = $row134

# ---------------------------------------------------------------------------
# 2) Append a new row 139 with a new match record.
#    Copy formatting (number formats / styles) from row 134 first, since
#    A139 and E139 need the same styles as the rest of column A / column E.
# ---------------------------------------------------------------------------
$ws.Range("A134").Copy()
$ws.Range("A139").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E134").Copy()
$ws.Range("E139").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Range("A139").Value = 138
$ws.Range("B139").Value = "bulgaria"
$ws.Range("C139").Value = "vtora-liga"
$ws.Range("D139").Value = "2023-2024"
$ws.Range("E139").Value = 45233.6875
$ws.Range("F139").Value = "Dunav Ruse"
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = "Chernomorets Balchik"
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 1.41
$ws.Range("K139").Value = "02/11/2023 04:42"
$ws.Range("L139").Value = 1.45
$ws.Range("M139").Value = "03/11/2023 16:28"
$ws.Range("N139").Value = 3.78
$ws.Range("O139").Value = "02/11/2023 04:42"
$ws.Range("P139").Value = 3.7
$ws.Range("Q139").Value = "03/11/2023 16:29"
$ws.Range("R139").Value = 5.63
$ws.Range("S139").Value = "02/11/2023 04:42"
$ws.Range("T139").Value = 6.84
$ws.Range("U139").Value = "03/11/2023 16:29"
$ws.Range("V139").Value = "https://www.betexplorer.com/football/bulgaria/vtora-liga/dunav-ruse-chernomorets-balchik/OdUxbUnf/"
